{"js": "// Word template \"ContentDemo.docx\" is being updated to demonstrate the new\n// Description field on the Graphic/Picture object. A new example placeholder\n// line - \"{txt:LocationPic.Description}\" followed by a manual line break -\n// is inserted right after the existing \"{pic:LocationPic:200}\" line (and its\n// line break) in the \"Picture\" demo paragraph, before the explanatory\n// \"{!:The pic:LocationPic:200 placeholder ...\" text.\n\nconst body = context.document.body;\n\n// Locate the existing placeholder text that marks our insertion point.\nconst results = body.search(\"{pic:LocationPic:200}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"{pic:LocationPic:200}\" in the document body.');\n}\n\nconst target = results.items[0];\n\n// Insert a manual line break followed by the new example placeholder text\n// immediately after the found text. The paragraph already has a line break\n// right after \"{pic:LocationPic:200}\" (before the red explanatory text), so\n// this new run ends up between that existing break and the explanatory text.\ntarget.insertText(\"\\u000b{txt:LocationPic.Description}\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word template \"ContentDemo.docx\" is being updated to demonstrate the new\n# Description field on the Graphic/Picture object. A new example placeholder\n# line - \"{txt:LocationPic.Description}\" - is inserted right after the\n# existing \"{pic:LocationPic:200}\" line in the \"Picture\" demo paragraph,\n# before the explanatory \"{!:The pic:LocationPic:200 placeholder ...\" text.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"{pic:LocationPic:200}\")\n\nif ($found) {\n    # Collapse to the end of the found text, then insert a manual line break\n    # (Chr 11, the same character Word uses for <w:br/>) followed by the new\n    # example placeholder text. The paragraph already has a line break right\n    # after \"{pic:LocationPic:200}\" (before the red explanatory text), so\n    # this ends up on its own line between that existing break and the\n    # explanatory text.\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter([char]11 + \"{txt:LocationPic.Description}\")\n}\n"}
